$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row after the existing data (column A)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1

$rows = @(
    @("lang_follow",   "Theo dõi",               "Follow"),
    @("lang_followed",  "Đã theo dõi",            "Followed"),
    @("lang_comment",   "Bình luận",              "Comment"),
    @("lang_search",    "Tìm kiếm địa chỉ…",      "Search")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
    $ws.Rows.Item($r).RowHeight = 14.9
}

# Select the final cell, matching the active selection state in the target file
$lastDataRow = $startRow + $rows.Count - 1
$ws.Range("C" + $lastDataRow).Select()
